$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell 2 4 '29.623.36'
Set-TextCell 2 5 '  -3.36%  '

Set-TextCell 3 4 '2.089.08'
Set-TextCell 3 5 '  -1.19%  '

Set-TextCell 5 4 '342.69'
Set-TextCell 5 5 '  -1.87%  '

Set-TextCell 6 4 '1.009'
Set-TextCell 6 5 '  -0.14%  '

Set-TextCell 7 4 '0.5152'
Set-TextCell 7 5 '  -2.07%  '

Set-TextCell 8 4 '0.4384'
Set-TextCell 8 5 '  -2.80%  '

Set-TextCell 9 2 'OKB'
Set-TextCell 9 3 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell 9 4 '52.28'
Set-TextCell 9 5 '  -2.55%  '

Set-TextCell 10 2 'Dogecoin'
Set-TextCell 10 3 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell 10 4 '0.09200'
Set-TextCell 10 5 '  +1.85%  '

Set-TextCell 11 4 '1.175'
Set-TextCell 11 5 '  +0.15%  '

Set-TextCell 12 4 '25.00'
Set-TextCell 12 5 '  +1.90%  '

Set-TextCell 13 4 '2.126.06'
Set-TextCell 13 5 '  +1.26%  '

Set-TextCell 14 4 '6.752'
Set-TextCell 14 5 '  -0.96%  '

Set-TextCell 15 4 '8.137'
Set-TextCell 15 5 '  +1.26%  '

Set-TextCell 16 5 '  -0.29%  '

Set-TextCell 17 4 '0.00001152'
Set-TextCell 17 5 '  -1.59%  '

Set-TextCell 18 5 '  -0.29%  '

Set-TextCell 19 4 '21.04'
Set-TextCell 19 5 '  +8.65%  '

Set-TextCell 20 4 '0.06641'
Set-TextCell 20 5 '  -0.94%  '

Set-TextCell 21 5 '  -0.21%  '

Set-TextCell 22 4 '6.173'
Set-TextCell 22 5 '  -1.94%  '

Set-TextCell 23 4 '29.669.66'
Set-TextCell 23 5 '  -3.42%  '

Set-TextCell 24 4 '12.60'
Set-TextCell 24 5 '  -2.02%  '

Set-TextCell 25 4 '2.313'
Set-TextCell 25 5 '  -3.16%  '

Set-TextCell 26 4 '2.345.12'
Set-TextCell 26 5 '  -0.41%  '

Set-TextCell 27 4 '21.82'
Set-TextCell 27 5 '  -2.22%  '

Set-TextCell 28 2 'Monero'
Set-TextCell 28 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 28 4 '162.43'
Set-TextCell 28 5 '  -1.83%  '

Set-TextCell 29 2 'LidoDAOToken'
Set-TextCell 29 3 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell 29 4 '2.525'
Set-TextCell 29 5 '  -0.30%  '

Set-TextCell 30 4 '132.47'
Set-TextCell 30 5 '  -2.31%  '

Set-TextCell 31 4 '1.128'
Set-TextCell 31 5 '  -5.20%  '

Set-TextCell 32 4 '0.1046'
Set-TextCell 32 5 '  -2.60%  '

Set-TextCell 33 4 '1.620'
Set-TextCell 33 5 '  -1.91%  '

Set-TextCell 34 4 '6.153'
Set-TextCell 34 5 '  -3.30%  '

Set-TextCell 35 4 '3.967'
Set-TextCell 35 5 '  -1.16%  '

Set-TextCell 36 4 '6.058'
Set-TextCell 36 5 '  +2.29%  '

Set-TextCell 37 4 '10.23'
Set-TextCell 37 5 '  +0.04%  '

Set-TextCell 38 4 '0.02564'
Set-TextCell 38 5 '  -3.27%  '

Set-TextCell 39 4 '0.06701'
Set-TextCell 39 5 '  -2.06%  '

Set-TextCell 40 4 '0.2236'
Set-TextCell 40 5 '  -3.31%  '

Set-TextCell 41 2 'Aptos'
Set-TextCell 41 3 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 41 4 '12.38'
Set-TextCell 41 5 '  -1.52%  '

Set-TextCell 42 2 'TheSandbox'
Set-TextCell 42 3 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell 42 4 '0.6856'
Set-TextCell 42 5 '  -0.54%  '

Set-TextCell 43 4 '1.287'
Set-TextCell 43 5 '  +0.51%  '

Set-TextCell 44 4 '0.6670'
Set-TextCell 44 5 '  +3.60%  '

Set-TextCell 45 4 '14.17'
Set-TextCell 45 5 '  -4.25%  '

Set-TextCell 46 4 '2.300'
Set-TextCell 46 5 '  -0.99%  '

Set-TextCell 47 5 '  -3.81%  '

Set-TextCell 48 4 '1.216'
Set-TextCell 48 5 '  -2.72%  '

Set-TextCell 49 4 '0.00000000335'
Set-TextCell 49 5 '  -5.17%  '

Set-TextCell 50 4 '81.35'
Set-TextCell 50 5 '  -1.21%  '

Set-TextCell 51 4 '1.167'
Set-TextCell 51 5 '  -2.12%  '
